$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new AQC overRide Report Download test-run rows (rows 125-140)
# mirrors the existing 8-row block pattern used for each run cycle,
# with an added "AQC overRide Report Download" step.
$data = @(
    ("Configuration setup", "Read all the configurations", "Pass", "22-01-2021 06:06:10 PM"),
    ("Login to Portal ", "Logged in sucessfully ", "Pass", "22-01-2021 06:06:38 PM"),
    ("Navigation to LATAM", "Successfully Navigated to LATAM ", "Pass", "22-01-2021 06:06:53 PM"),
    ("Search of created WO ", "Search is Working", "Pass", "22-01-2021 06:07:06 PM"),
    ("Segment level Ingest", "Segment level Ingest is completed", "Pass", "22-01-2021 06:07:34 PM"),
    ("AQC overRide is clicked ", "AQC over Ride performed ", "Pass", "22-01-2021 06:07:52 PM"),
    ("AQC overRide Report Download ", "AQC overRide Report Download ", "Pass", "22-01-2021 06:07:59 PM"),
    ("AQC overRide pass is enabled ", "AQC over ride Pass", "Fail", "22-01-2021 06:07:59 PM"),
    ("Configuration setup", "Read all the configurations", "Pass", "22-01-2021 06:26:11 PM"),
    ("Login to Portal ", "Logged in sucessfully ", "Pass", "22-01-2021 06:26:37 PM"),
    ("Navigation to LATAM", "Successfully Navigated to LATAM ", "Pass", "22-01-2021 06:26:51 PM"),
    ("Search of created WO ", "Search is Working", "Pass", "22-01-2021 06:27:04 PM"),
    ("Segment level Ingest", "Segment level Ingest is completed", "Pass", "22-01-2021 06:27:32 PM"),
    ("AQC overRide is clicked ", "AQC over Ride performed ", "Pass", "22-01-2021 06:27:50 PM"),
    ("AQC overRide Report Download ", "AQC overRide Report Download ", "Pass", "22-01-2021 06:27:57 PM"),
    ("AQC overRide pass is enabled ", "AQC over ride Pass", "Fail", "22-01-2021 06:27:57 PM")
)

$startRow = 125
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Range("A$row").Value = $data[$i][0]
    $ws.Range("B$row").Value = $data[$i][1]
    $ws.Range("C$row").Value = $data[$i][2]
    $ws.Range("D$row").Value = $data[$i][3]
}
